$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = 93
$ws.Range("D21").Value = 86
$ws.Range("E21").Value = 7
$ws.Range("F21").Value = 24.64183381088825
